$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("L2").Value = 4851
$ws.Range("L3").Value = 5205
$ws.Range("D4").Value = 1996
$ws.Range("L4").Value = 1274
$ws.Range("L5").Value = 306
$ws.Range("L6").Value = 4412
$ws.Range("D7").Value = 28187
$ws.Range("L7").Value = 16048

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("L2").Value = 55
$ws.Range("L7").Value = 182

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("L2").Value = 311
$ws.Range("L3").Value = 361
$ws.Range("L7").Value = 1068

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("L2").Value = 111
$ws.Range("L3").Value = 143
$ws.Range("L6").Value = 83
$ws.Range("L7").Value = 359

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("L2").Value = 201
$ws.Range("L3").Value = 252
$ws.Range("L7").Value = 734

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("L6").Value = 48
$ws.Range("L7").Value = 221

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("L2").Value = 181
$ws.Range("L6").Value = 169
$ws.Range("L7").Value = 606

$ws = $wb.Worksheets.Item('New City')
$ws.Range("L2").Value = 112
$ws.Range("L7").Value = 313

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("L4").Value = 61
$ws.Range("L6").Value = 122
$ws.Range("L7").Value = 526
$ws.Range("L8").Value = 1068
$ws.Range("L9").Value = 93
$ws.Range("L11").Value = 260
$ws.Range("L14").Value = 86
$ws.Range("L15").Value = 118
$ws.Range("L18").Value = 113
$ws.Range("L19").Value = 440
$ws.Range("L29").Value = 877
$ws.Range("L33").Value = 734
$ws.Range("L36").Value = 209
$ws.Range("L37").Value = 606
$ws.Range("L42").Value = 521
$ws.Range("L43").Value = 117
$ws.Range("L45").Value = 30
$ws.Range("L48").Value = 208
$ws.Range("L51").Value = 203
$ws.Range("L52").Value = 324
$ws.Range("L53").Value = 182
$ws.Range("L54").Value = 336
$ws.Range("L56").Value = 15
$ws.Range("D63").Value = 376
$ws.Range("L63").Value = 47
$ws.Range("L65").Value = 313
$ws.Range("L67").Value = 553
$ws.Range("L79").Value = 422
$ws.Range("L83").Value = 359
$ws.Range("L85").Value = 824
$ws.Range("L89").Value = 232
$ws.Range("L91").Value = 221
$ws.Range("L95").Value = 221
$ws.Range("L96").Value = 180
$ws.Range("L97").Value = 137
$ws.Range("L98").Value = 87
$ws.Range("D101").Value = 28187
$ws.Range("L101").Value = 16048

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("L2").Value = 162
$ws.Range("L6").Value = 127
$ws.Range("L7").Value = 553

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("L3").Value = 82
$ws.Range("L7").Value = 336

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("L2").Value = 263
$ws.Range("L5").Value = 15
$ws.Range("L7").Value = 877

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("L6").Value = 87
$ws.Range("L7").Value = 208

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("L2").Value = 156
$ws.Range("L3").Value = 137
$ws.Range("L7").Value = 440

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("L6").Value = 22
$ws.Range("L7").Value = 86

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("L3").Value = 35
$ws.Range("L7").Value = 122

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("L2").Value = 151
$ws.Range("L3").Value = 173
$ws.Range("L7").Value = 521

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("L3").Value = 53
$ws.Range("L7").Value = 180

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("L2").Value = 78
$ws.Range("L7").Value = 221

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("L6").Value = 92
$ws.Range("L7").Value = 422

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("L6").Value = 20
$ws.Range("L7").Value = 113

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("L2").Value = 77
$ws.Range("L7").Value = 209

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("L2").Value = 180
$ws.Range("L7").Value = 526

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("L3").Value = 39
$ws.Range("L7").Value = 118

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("L2").Value = 21
$ws.Range("L7").Value = 87

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("L2").Value = 98
$ws.Range("L7").Value = 260

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range("L2").Value = 26
$ws.Range("L7").Value = 93

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("L2").Value = 32
$ws.Range("L7").Value = 137

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("L4").Value = 37
$ws.Range("L7").Value = 232

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("L2").Value = 60
$ws.Range("L7").Value = 203

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("L3").Value = 36
$ws.Range("L7").Value = 117

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("L6").Value = 175
$ws.Range("L7").Value = 824

$ws = $wb.Worksheets.Item('Jackson Park')
$ws.Range("L6").Value = 7
$ws.Range("L7").Value = 30

$ws = $wb.Worksheets.Item('Magnificent Mile')
$ws.Range("L6").Value = 8
$ws.Range("L7").Value = 15

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("L6").Value = 89
$ws.Range("L7").Value = 324

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Range("L2").Value = 20
$ws.Range("L7").Value = 61
